{"js": "// Load the first two paragraphs of the document body.\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst titlePara = paras.items[0];   // \"ME3023\" + tabs + \"Electronics Basics Background\"\nconst infoPara = paras.items[1];    // \"Information: \" (bold)\n\n// 0) Protect the \"Electronics Basics\" / \" Background\" run boundary with a\n//    throwaway bookmark so later edits to the title paragraph don't coalesce\n//    the two runs into one.\nconst ebEnd = titlePara.search(\"Electronics Basics\", { matchCase: true });\nawait context.sync();\nebEnd.items[0].getRange(\"End\").insertBookmark(\"TEMP_RUN_ANCHOR\");\nawait context.sync();\n\n// 1) Remove the \"ME3023\" text from the title paragraph.\nconst mePrefix = titlePara.search(\"ME3023\", { matchCase: true });\nawait context.sync();\nmePrefix.items[0].delete();\nawait context.sync();\n\n// 2) Remove the remaining tab characters in the title paragraph\n//    (delete from the end backwards so earlier matches stay valid).\nconst tabs = titlePara.search(\"\\t\", { matchCase: true });\ntabs.load(\"items\");\nawait context.sync();\nfor (let i = tabs.items.length - 1; i >= 0; i--) {\n  tabs.items[i].delete();\n  await context.sync();\n}\n\n// 3) Drop the old \"_GoBack\" bookmark that used to sit in the title paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3b) Remove the temporary anchor bookmark now that the runs are split.\ncontext.document.deleteBookmark(\"TEMP_RUN_ANCHOR\");\nawait context.sync();\n\n// 4) Center the title paragraph (\"Electronics Basics Background\").\ntitlePara.alignment = Word.Alignment.centered;\nawait context.sync();\n\n// 5) Remove the trailing space that followed \"Information: \".\nconst trailingSpace = infoPara.search(\" \", { matchCase: true });\ntrailingSpace.load(\"items\");\nawait context.sync();\nfor (let i = trailingSpace.items.length - 1; i >= 0; i--) {\n  trailingSpace.items[i].delete();\n  await context.sync();\n}\n\n// 6) Re-insert the \"_GoBack\" bookmark right after \"Information:\" (end of\n//    the now-trimmed second paragraph).\nconst endOfInfo = infoPara.getRange(\"End\");\nendOfInfo.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Title paragraph: \"ME3023\" + tabs + \"Electronics Basics Background\" ---\n$titlePara = $d.Paragraphs(1).Range\n\n# 0) Protect the \"Electronics Basics\" / \" Background\" run boundary with a\n#    throwaway bookmark so later edits don't coalesce the two runs into one.\n$anchorRng = $titlePara.Duplicate\n$anchorRng.Find.Execute(\"Electronics Basics\") | Out-Null\n$anchorRng.Collapse(0)\n$d.Bookmarks.Add(\"TEMP_RUN_ANCHOR\", $anchorRng) | Out-Null\n\n# 1) Remove the \"ME3023\" text (scoped to this paragraph only - \"ME3023\"\n#    also appears later in the document and must stay untouched there).\n$meRng = $titlePara.Duplicate\n$meRng.Find.Execute(\"ME3023\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\n# 2) Remove the remaining tab characters in the title paragraph (scoped -\n#    other paragraphs elsewhere also contain tabs that must stay untouched).\n$tabRng = $titlePara.Duplicate\n$tabRng.Find.Execute(\"^t\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\n# 3) Drop the old \"_GoBack\" bookmark that used to sit in the title paragraph.\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# 3b) Remove the temporary anchor bookmark now that the runs are split.\n$d.Bookmarks(\"TEMP_RUN_ANCHOR\").Delete()\n\n# 4) Center the title paragraph (\"Electronics Basics Background\").\n$d.Paragraphs(1).Alignment = 1\n\n# --- Second paragraph: bold \"Information: \" ---\n$infoPara = $d.Paragraphs(2).Range\n\n# 5) Remove the trailing space after \"Information:\" and relocate the\n#    \"_GoBack\" bookmark there. A temporary \"X\" marker is used so the\n#    collapsed range sits strictly inside the paragraph (not at the very\n#    end, right before the paragraph mark) while the bookmark is created;\n#    it is deleted immediately afterwards.\n$infoRng = $infoPara.Duplicate\n$infoRng.Find.Execute(\"Information: \", $false, $false, $false, $false, $false, $true, 1, $false, \"Information:X\", 2) | Out-Null\n\n$markRng = $d.Paragraphs(2).Range.Duplicate\n$markRng.Find.Execute(\"Information:\") | Out-Null\n$markRng.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $markRng) | Out-Null\n\n$cleanupRng = $d.Paragraphs(2).Range.Duplicate\n$cleanupRng.Find.Execute(\"X\", $true, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\nWrite-Output \"done\"\n"}
